$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9791771173477173
$ws.Range("B1").Value = 3.406835079193115
$ws.Range("C1").Value = 4.397721767425537
$ws.Range("D1").Value = 1.920611977577209
$ws.Range("E1").Value = 0.7703719139099121
